$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 changes
$ws.Range("D2").Value = 0.008650000000000001
$ws.Range("E2").Value = 0.2570000000000001
$ws.Range("G2").Value = 0.2540117994100295
$ws.Range("H2").Value = 0.2540117994100295
$ws.Range("I2").Value = 0.07415929203539824
$ws.Range("J2").Value = 0.05833450586884883
$ws.Range("K2").Value = 2.013
$ws.Range("L2").Value = 0.05938053097345133
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0
$ws.Range("T2").ClearContents()
$ws.Range("U2").Value = 11.41
$ws.Range("V2").Value = 0.161614730878187
$ws.Range("W2").Value = 0.061272412285128
$ws.Range("X2").Value = 0.08364952059172212
$ws.Range("Y2").Value = -0.02237710830659412
$ws.Range("Z2").Value = 1.233804047168438
$ws.Range("AA2").Value = 0.07464194047739327
$ws.Range("AB2").Value = 0.06758500954196445
$ws.Range("AC2").Value = 0.007056930935428824
$ws.Range("AD2").Value = 4.789
$ws.Range("AF2").Value = 4.789
$ws.Range("AG2").Value = -6.621
$ws.Range("AH2").Value = 0.06352385626550291
$ws.Range("AI2").Value = 0.08808328275303941
$ws.Range("AJ2").Value = -0.1034870816986824
$ws.Range("AK2").Value = -0.1541236993412324
$ws.Range("AL2").Value = 0.149
$ws.Range("AM2").Value = 0.149
$ws.Range("AN2").Value = 1.6054307743882
$ws.Range("AO2").Value = 16.87248322147651
$ws.Range("AP2").Value = -2.219577606436474
$ws.Range("AQ2").Value = 16.87248322147651

# Row 3 changes
$ws.Range("B3").Value = "Minerva Insurance Company Public Ltd (CSE:MINE)"
$ws.Range("D3").Value = 0.0558
$ws.Range("E3").Value = 0.7490000000000001
$ws.Range("G3").Value = 0.0138125
$ws.Range("H3").Value = 0.0138125
$ws.Range("I3").Value = 0.030875
$ws.Range("J3").Value = 0.030875
$ws.Range("K3").Value = 0.653
$ws.Range("L3").Value = 0.0408125
$ws.Range("M3").Value = -0
$ws.Range("N3").Value = -0
$ws.Range("O3").Value = -0
$ws.Range("P3").Value = -0
$ws.Range("Q3").Value = -0
$ws.Range("R3").Value = -0
$ws.Range("T3").ClearContents()
$ws.Range("U3").Value = 4.5
$ws.Range("V3").Value = 1.022727272727273
$ws.Range("W3").Value = 0.08896457765667576
$ws.Range("X3").Value = 0.1075599347608003
$ws.Range("Y3").Value = -0.01859535710412458
$ws.Range("Z3").Value = 3.16205533596838
$ws.Range("AA3").Value = 0.09762845849802373
$ws.Range("AB3").Value = 0.07543656774464441
$ws.Range("AC3").Value = 0.02219189075337932
$ws.Range("AD3").Value = 4.77
$ws.Range("AF3").Value = 4.77
$ws.Range("AG3").Value = 0.2699999999999996
$ws.Range("AH3").Value = 0.520174482006543
$ws.Range("AI3").Value = 0.3712062256809338
$ws.Range("AJ3").Value = 0.05781584582441104
$ws.Range("AK3").Value = 0.03233532934131732
$ws.Range("AL3").Value = 0.149
$ws.Range("AM3").Value = 0.149
$ws.Range("AN3").Value = 6.170763260025873
$ws.Range("AO3").Value = 3.315436241610739
$ws.Range("AP3").Value = 0.3492884864165583
$ws.Range("AQ3").Value = 3.315436241610739

# Row 4 changes
$ws.Range("B4").Value = "Atlantic Insurance Company Public Limited (CSE:ATL)"
$ws.Range("D4").Value = -0.0385
$ws.Range("E4").Value = -0.235
$ws.Range("G4").Value = 0.4687150837988828
$ws.Range("H4").Value = 0.4687150837988828
$ws.Range("I4").Value = 0.1128491620111732
$ws.Range("J4").Value = 0.06468759496037961
$ws.Range("K4").Value = 1.36
$ws.Range("L4").Value = 0.07597765363128493
$ws.Range("O4").Value = -0
$ws.Range("R4").Value = -0
$ws.Range("U4").Value = 6.91
$ws.Range("V4").Value = 0.104380664652568
$ws.Range("W4").Value = 0.03358024691358025
$ws.Range("X4").Value = 0.05973910642264389
$ws.Range("Y4").Value = -0.02615885950906364
$ws.Range("Z4").Value = 0.7985367594575303
$ws.Range("AA4").Value = 0.05165542245676281
$ws.Range("AB4").Value = 0.05973345133928448
$ws.Range("AC4").Value = -0.008078028882521675
$ws.Range("AD4").Value = 0.019
$ws.Range("AF4").Value = 0.019
$ws.Range("AG4").Value = -6.891
$ws.Range("AH4").Value = 0.0002869267128769688
$ws.Range("AI4").Value = 0.0004576218117006671
$ws.Range("AJ4").Value = -0.1161880996138866
$ws.Range("AK4").Value = -0.1991100580773787
$ws.Range("AL4").Value = 0
$ws.Range("AM4").Value = 0
$ws.Range("AN4").Value = 0.008597285067873304
$ws.Range("AO4").ClearContents()
$ws.Range("AP4").Value = -3.118099547511312
$ws.Range("AQ4").ClearContents()

